$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($cell, [string]$val)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextCell $ws.Cells.Item(2, 4) '26.866.07'
Set-TextCell $ws.Cells.Item(2, 5) '  +1.12%  '
Set-TextCell $ws.Cells.Item(3, 4) '1.840.72'
Set-TextCell $ws.Cells.Item(3, 5) '  +1.51%  '
Set-TextCell $ws.Cells.Item(4, 4) '1.007'
Set-TextCell $ws.Cells.Item(4, 5) '  +0.05%  '
Set-TextCell $ws.Cells.Item(5, 4) '308.87'
Set-TextCell $ws.Cells.Item(5, 5) '  +1.27%  '
Set-TextCell $ws.Cells.Item(6, 5) '  +0.03%  '
Set-TextCell $ws.Cells.Item(7, 4) '0.4699'
Set-TextCell $ws.Cells.Item(7, 5) '  +0.85%  '
Set-TextCell $ws.Cells.Item(8, 4) '0.3660'
Set-TextCell $ws.Cells.Item(8, 5) '  +2.47%  '
Set-TextCell $ws.Cells.Item(9, 4) '0.07168'
Set-TextCell $ws.Cells.Item(9, 5) '  +0.80%  '
Set-TextCell $ws.Cells.Item(10, 4) '0.9179'
Set-TextCell $ws.Cells.Item(10, 5) '  +2.02%  '
Set-TextCell $ws.Cells.Item(11, 4) '1.927.31'
Set-TextCell $ws.Cells.Item(11, 5) '  +5.28%  '
Set-TextCell $ws.Cells.Item(12, 4) '19.52'
Set-TextCell $ws.Cells.Item(12, 5) '  +0.93%  '
Set-TextCell $ws.Cells.Item(13, 4) '0.07607'
Set-TextCell $ws.Cells.Item(13, 5) '  -1.84%  '
Set-TextCell $ws.Cells.Item(14, 4) '5.279'
Set-TextCell $ws.Cells.Item(14, 5) '  +0.75%  '
Set-TextCell $ws.Cells.Item(15, 4) '6.404'
Set-TextCell $ws.Cells.Item(15, 5) '  +1.95%  '
Set-TextCell $ws.Cells.Item(16, 4) '87.95'
Set-TextCell $ws.Cells.Item(16, 5) '  +0.49%  '
Set-TextCell $ws.Cells.Item(17, 5) '  +0.10%  '
Set-TextCell $ws.Cells.Item(18, 5) '  +0.96%  '
Set-TextCell $ws.Cells.Item(19, 5) '  +0.02%  '
Set-TextCell $ws.Cells.Item(20, 4) '26.889.21'
Set-TextCell $ws.Cells.Item(20, 5) '  +1.08%  '
Set-TextCell $ws.Cells.Item(21, 4) '14.52'
Set-TextCell $ws.Cells.Item(21, 5) '  +2.89%  '
Set-TextCell $ws.Cells.Item(22, 4) '5.002'
Set-TextCell $ws.Cells.Item(22, 5) '  +0.54%  '
Set-TextCell $ws.Cells.Item(23, 4) '10.60'
Set-TextCell $ws.Cells.Item(23, 5) '  +0.51%  '
Set-TextCell $ws.Cells.Item(24, 4) '1.923'
Set-TextCell $ws.Cells.Item(24, 5) '  +0.11%  '
Set-TextCell $ws.Cells.Item(25, 4) '151.75'
Set-TextCell $ws.Cells.Item(25, 5) '  -0.13%  '
Set-TextCell $ws.Cells.Item(26, 4) '18.14'
Set-TextCell $ws.Cells.Item(26, 5) '  +1.61%  '
Set-TextCell $ws.Cells.Item(27, 4) '2.001'
Set-TextCell $ws.Cells.Item(27, 5) '  -0.49%  '
Set-TextCell $ws.Cells.Item(28, 4) '114.09'
Set-TextCell $ws.Cells.Item(28, 5) '  +1.31%  '
Set-TextCell $ws.Cells.Item(29, 4) '4.848'
Set-TextCell $ws.Cells.Item(29, 5) '  +0.94%  '
Set-TextCell $ws.Cells.Item(30, 4) '0.08820'
Set-TextCell $ws.Cells.Item(31, 4) '3.264'
Set-TextCell $ws.Cells.Item(31, 5) '  +4.90%  '
Set-TextCell $ws.Cells.Item(32, 2) 'RenderToken'
Set-TextCell $ws.Cells.Item(32, 3) 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextCell $ws.Cells.Item(32, 4) '2.801'
Set-TextCell $ws.Cells.Item(32, 5) '  +3.89%  '
Set-TextCell $ws.Cells.Item(33, 2) 'ARBITRUM'
Set-TextCell $ws.Cells.Item(33, 3) 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextCell $ws.Cells.Item(33, 4) '1.167'
Set-TextCell $ws.Cells.Item(33, 5) '  +4.41%  '
Set-TextCell $ws.Cells.Item(34, 2) 'ImmutableX'
Set-TextCell $ws.Cells.Item(34, 3) 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextCell $ws.Cells.Item(34, 4) '0.7417'
Set-TextCell $ws.Cells.Item(34, 5) '  +0.81%  '
Set-TextCell $ws.Cells.Item(35, 2) 'Filecoin'
Set-TextCell $ws.Cells.Item(35, 3) 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextCell $ws.Cells.Item(35, 4) '4.478'
Set-TextCell $ws.Cells.Item(35, 5) '  +1.09%  '
Set-TextCell $ws.Cells.Item(36, 4) '1.085'
Set-TextCell $ws.Cells.Item(36, 5) '  +1.29%  '
Set-TextCell $ws.Cells.Item(37, 4) '0.05267'
Set-TextCell $ws.Cells.Item(37, 5) '  +4.03%  '
Set-TextCell $ws.Cells.Item(38, 2) 'VeChain'
Set-TextCell $ws.Cells.Item(38, 3) 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextCell $ws.Cells.Item(38, 4) '0.01943'
Set-TextCell $ws.Cells.Item(38, 5) '  +0.89%  '
Set-TextCell $ws.Cells.Item(39, 2) 'MXToken'
Set-TextCell $ws.Cells.Item(39, 3) 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextCell $ws.Cells.Item(39, 4) '2.969'
Set-TextCell $ws.Cells.Item(39, 5) '  +2.01%  '
Set-TextCell $ws.Cells.Item(40, 4) '0.5178'
Set-TextCell $ws.Cells.Item(40, 5) '  +2.67%  '
Set-TextCell $ws.Cells.Item(41, 4) '6.880'
Set-TextCell $ws.Cells.Item(41, 5) '  +0.94%  '
Set-TextCell $ws.Cells.Item(42, 4) '0.1511'
Set-TextCell $ws.Cells.Item(42, 5) '  +1.38%  '
Set-TextCell $ws.Cells.Item(43, 4) '8.147'
Set-TextCell $ws.Cells.Item(43, 5) '  +2.28%  '
Set-TextCell $ws.Cells.Item(44, 4) '10.46'
Set-TextCell $ws.Cells.Item(44, 5) '  +5.18%  '
Set-TextCell $ws.Cells.Item(45, 4) '0.4679'
Set-TextCell $ws.Cells.Item(45, 5) '  +0.33%  '
Set-TextCell $ws.Cells.Item(46, 4) '1.008'
Set-TextCell $ws.Cells.Item(46, 5) '  +0.07%  '
Set-TextCell $ws.Cells.Item(47, 4) '101.56'
Set-TextCell $ws.Cells.Item(47, 5) '  +2.82%  '
Set-TextCell $ws.Cells.Item(48, 4) '1.593'
Set-TextCell $ws.Cells.Item(48, 5) '  +1.93%  '
Set-TextCell $ws.Cells.Item(49, 4) '65.31'
Set-TextCell $ws.Cells.Item(49, 5) '  +2.70%  '
Set-TextCell $ws.Cells.Item(50, 5) '  +0.21%  '
Set-TextCell $ws.Cells.Item(51, 4) '0.8840'
Set-TextCell $ws.Cells.Item(51, 5) '  +4.24%  '
